$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 45992
$ws.Range("A25").NumberFormat = $ws.Range("A24").NumberFormat

$ws.Range("B25").Value = 6435
$ws.Range("C25").Value = 1003
$ws.Range("D25").Value = 5994969
$ws.Range("E25").Value = 931.6191142191142
$ws.Range("F25").Value = 9.234425394669831
$ws.Range("G25").Value = 6.929637526652455
$ws.Range("H25").Value = 25.55313825828101
